$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Helper to write a value as plain text (avoiding automatic number/date
# inference), then strip the explicit style that NumberFormat leaves behind
# so the cell stays visually/structurally identical to its neighbours.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Grow the table by two rows (12 and 13), keeping the existing table
# column/autofilter definitions untouched.
[void]$lo.ListRows.Add()
[void]$lo.ListRows.Add()

# Row 12: Maçã
Set-TextValue $ws.Cells.Item(12, 1) "Maçã"
Set-TextValue $ws.Cells.Item(12, 2) "001"
Set-TextValue $ws.Cells.Item(12, 3) "17/05/2025"
Set-TextValue $ws.Cells.Item(12, 4) "super"
Set-TextValue $ws.Cells.Item(12, 5) "Alimento"
Set-TextValue $ws.Cells.Item(12, 6) "2"
Set-TextValue $ws.Cells.Item(12, 7) "Verde`n"

# Row 13: Melão
Set-TextValue $ws.Cells.Item(13, 1) "Melão"
Set-TextValue $ws.Cells.Item(13, 2) "777"
Set-TextValue $ws.Cells.Item(13, 3) "02/06/2025"
Set-TextValue $ws.Cells.Item(13, 4) "super"
Set-TextValue $ws.Cells.Item(13, 5) "Alimento"
Set-TextValue $ws.Cells.Item(13, 6) "2"
Set-TextValue $ws.Cells.Item(13, 7) "`n"

# Writing multi-line text triggers an automatic row-height recalculation;
# re-running AutoFit puts the (wrap-off) rows back to the standard height
# so no spurious custom row height is left behind.
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()
